# Add a new "2022-Q1" sheet with fund holdings data, positioned before "总计",
# and insert a new leading row "2022-Q1" into the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# Locate the existing "总计" sheet (currently the last sheet).
$totalSheet = $wb.Worksheets.Item("总计")

# Create the new sheet right before "总计" (Add's first arg is the "Before" sheet).
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: sheet handles obtained before an insert/add can re-seat to a
# different sheet once the collection shifts, so re-fetch "总计" by name
# now that the new sheet exists.
$totalSheet = $wb.Worksheets.Item("总计")

# --- Populate "2022-Q1" sheet (same layout as the other quarterly sheets) ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'160325"
$newSheet.Range("C2").Value = "华夏创业板两年定期开放混合"
$newSheet.Range("D2").Value = "'27.39"
$newSheet.Range("E2").Value = "'90.77"
$newSheet.Range("F2").Value = "'2.19"
$newSheet.Range("G2").Value = "'0.5998"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'003842"
$newSheet.Range("C3").Value = "中邮景泰灵活配置混合A"
$newSheet.Range("D3").Value = "'5.76"
$newSheet.Range("E3").Value = "'33.95"
$newSheet.Range("F3").Value = "'1.21"
$newSheet.Range("G3").Value = "'0.0697"
$newSheet.Range("H3").Value = 5

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'003843"
$newSheet.Range("C4").Value = "中邮景泰灵活配置混合C"
$newSheet.Range("D4").Value = "'0.37"
$newSheet.Range("E4").Value = "'33.95"
$newSheet.Range("F4").Value = "'1.21"
$newSheet.Range("G4").Value = "'0.0045"
$newSheet.Range("H4").Value = 5

# --- Normalize formatting on "2022-Q1" to match the other quarterly sheets:
#     bordered/bold style on the header row + index column, no special
#     (e.g. quote-prefix) style anywhere else.
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

$totalSheet.Range("A1").Copy()
$newSheet.Range("B2:G4").PasteSpecial(-4122)
$newSheet.Range("H2:H4").PasteSpecial(-4122)

# --- Update "总计" sheet: insert a new row 2 for "2022-Q1" ---
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.67

# Restore the bordered/bold style on the new index cell (A2), matching the
# other rows in this column.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Renumber the index column (A) for the rows pushed down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
